$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.180921052631579
$ws.Range("C2").Value = 0.5789473684210527
$ws.Range("J2").Value = 0.02302631578947368
$ws.Range("O2").Value = 0.003289473684210526
$ws.Range("P2").Value = 0.1414473684210526
$ws.Range("S2").Value = 0.07236842105263158
$ws.Range("B3").Value = 0.01666666666666667
$ws.Range("C3").Value = 0.01666666666666667
$ws.Range("J3").Value = 0.06666666666666667
$ws.Range("P3").Value = 0.7333333333333333
$ws.Range("S3").Value = 0.1666666666666667
$ws.Range("J4").Value = 0.0576923076923077
$ws.Range("P4").Value = 0.5769230769230769
$ws.Range("S4").Value = 0.3653846153846154
$ws.Range("B6").Value = 0.07352941176470588
$ws.Range("D6").Value = 0.02205882352941177
$ws.Range("F6").Value = 0.09558823529411764
$ws.Range("J6").Value = 0.2536764705882353
$ws.Range("O6").Value = 0.02573529411764706
$ws.Range("Q6").Value = 0.1397058823529412
$ws.Range("R6").Value = 0.0625
$ws.Range("S6").Value = 0.3272058823529412
$ws.Range("B7").Value = 0.08733624454148471
$ws.Range("D7").Value = 0.02620087336244541
$ws.Range("F7").Value = 0.05240174672489083
$ws.Range("J7").Value = 0.1135371179039301
$ws.Range("O7").Value = 0.04366812227074236
$ws.Range("Q7").Value = 0.222707423580786
$ws.Range("R7").Value = 0.07423580786026202
$ws.Range("S7").Value = 0.3799126637554585
$ws.Range("B8").Value = 0.0954356846473029
$ws.Range("D8").Value = 0.01037344398340249
$ws.Range("F8").Value = 0.05186721991701245
$ws.Range("J8").Value = 0.1348547717842324
$ws.Range("O8").Value = 0.01037344398340249
$ws.Range("Q8").Value = 0.1742738589211618
$ws.Range("R8").Value = 0.1286307053941909
$ws.Range("S8").Value = 0.3941908713692946
$ws.Range("B9").Value = 0.04102564102564103
$ws.Range("D9").Value = 0.03589743589743589
$ws.Range("F9").Value = 0.08717948717948718
$ws.Range("J9").Value = 0.1384615384615385
$ws.Range("O9").Value = 0.02051282051282051
$ws.Range("Q9").Value = 0.2051282051282051
$ws.Range("R9").Value = 0.09743589743589744
$ws.Range("S9").Value = 0.3743589743589744
$ws.Range("B10").Value = 0.1014492753623188
$ws.Range("D10").Value = 0.02208419599723948
$ws.Range("E10").Value = 0.001380262249827467
$ws.Range("F10").Value = 0.07246376811594203
$ws.Range("J10").Value = 0.1331953071083506
$ws.Range("O10").Value = 0.02277432712215321
$ws.Range("Q10").Value = 0.2111801242236025
$ws.Range("R10").Value = 0.09109730848861283
$ws.Range("S10").Value = 0.3443754313319531
$ws.Range("F11").Value = 0.002857142857142857
$ws.Range("G11").Value = 0.1571428571428571
$ws.Range("J11").Value = 0.1028571428571429
$ws.Range("K11").Value = 0.2057142857142857
$ws.Range("L11").Value = 0.5228571428571429
$ws.Range("S11").Value = 0.008571428571428572
$ws.Range("G12").Value = 0.7202072538860104
$ws.Range("J12").Value = 0.1968911917098446
$ws.Range("K12").Value = 0.0155440414507772
$ws.Range("L12").Value = 0.0310880829015544
$ws.Range("S12").Value = 0.03626943005181347
$ws.Range("G13").Value = 0.7586206896551724
$ws.Range("J13").Value = 0.1551724137931035
$ws.Range("S13").Value = 0.08620689655172414
$ws.Range("F15").Value = 0.03076923076923077
$ws.Range("H15").Value = 0.1269230769230769
$ws.Range("I15").Value = 0.07692307692307693
$ws.Range("J15").Value = 0.3653846153846154
$ws.Range("K15").Value = 0.09615384615384616
$ws.Range("M15").Value = 0.003846153846153846
$ws.Range("O15").Value = 0.07307692307692308
$ws.Range("S15").Value = 0.2269230769230769
$ws.Range("F16").Value = 0.03015075376884422
$ws.Range("H16").Value = 0.2261306532663317
$ws.Range("I16").Value = 0.06532663316582915
$ws.Range("J16").Value = 0.3266331658291458
$ws.Range("K16").Value = 0.1306532663316583
$ws.Range("M16").Value = 0.02512562814070352
$ws.Range("O16").Value = 0.05025125628140704
$ws.Range("S16").Value = 0.1457286432160804
$ws.Range("F17").Value = 0.01740812379110251
$ws.Range("H17").Value = 0.1624758220502901
$ws.Range("I17").Value = 0.05415860735009671
$ws.Range("J17").Value = 0.4506769825918762
$ws.Range("K17").Value = 0.1025145067698259
$ws.Range("M17").Value = 0.02901353965183752
$ws.Range("O17").Value = 0.08123791102514506
$ws.Range("S17").Value = 0.1025145067698259
$ws.Range("F18").Value = 0.03643724696356275
$ws.Range("H18").Value = 0.1700404858299595
$ws.Range("I18").Value = 0.0931174089068826
$ws.Range("J18").Value = 0.4089068825910931
$ws.Range("K18").Value = 0.1012145748987854
$ws.Range("M18").Value = 0.02834008097165992
$ws.Range("O18").Value = 0.06072874493927125
$ws.Range("S18").Value = 0.1012145748987854
$ws.Range("F19").Value = 0.02061068702290076
$ws.Range("H19").Value = 0.2137404580152672
$ws.Range("I19").Value = 0.083206106870229
$ws.Range("J19").Value = 0.3732824427480916
$ws.Range("K19").Value = 0.1083969465648855
$ws.Range("M19").Value = 0.02519083969465649
$ws.Range("N19").Value = 0.001526717557251908
$ws.Range("O19").Value = 0.0633587786259542
$ws.Range("S19").Value = 0.1106870229007634
